$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

Set-TextValue $ws.Cells.Item(2, 4) "21.999.97"
$ws.Cells.Item(2, 5).Value = "  -1.61%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.550.19"
$ws.Cells.Item(3, 5).Value = "  -1.02%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 5).Value = "  -0.21%  "
Set-TextValue $ws.Cells.Item(6, 4) "287.77"
$ws.Cells.Item(6, 5).Value = "  +0.22%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.3917"
$ws.Cells.Item(7, 5).Value = "  +3.71%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3195"
$ws.Cells.Item(8, 5).Value = "  -2.30%  "
Set-TextValue $ws.Cells.Item(9, 4) "42.02"
$ws.Cells.Item(9, 5).Value = "  -7.63%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.07250"
$ws.Cells.Item(10, 5).Value = "  -2.28%  "
$ws.Cells.Item(11, 5).Value = "  -4.23%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.002"
$ws.Cells.Item(12, 5).Value = "  +0.07%  "
$ws.Cells.Item(13, 5).Value = "  -7.60%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.603"
$ws.Cells.Item(14, 5).Value = "  -4.40%  "
Set-TextValue $ws.Cells.Item(15, 4) "6.633"
$ws.Cells.Item(15, 5).Value = "  -2.58%  "
Set-TextValue $ws.Cells.Item(16, 4) "0.00001120"
$ws.Cells.Item(16, 5).Value = "  +2.27%  "
Set-TextValue $ws.Cells.Item(17, 4) "1.549.03"
$ws.Cells.Item(17, 5).Value = "  -0.15%  "
Set-TextValue $ws.Cells.Item(18, 4) "0.06578"
$ws.Cells.Item(18, 5).Value = "  -2.22%  "
Set-TextValue $ws.Cells.Item(19, 4) "83.42"
$ws.Cells.Item(19, 5).Value = "  -3.01%  "
$ws.Cells.Item(20, 5).Value = "  -0.31%  "
Set-TextValue $ws.Cells.Item(21, 4) "6.272"
$ws.Cells.Item(21, 5).Value = "  -1.52%  "
Set-TextValue $ws.Cells.Item(22, 4) "15.67"
$ws.Cells.Item(22, 5).Value = "  -3.75%  "
Set-TextValue $ws.Cells.Item(23, 4) "11.20"
$ws.Cells.Item(23, 5).Value = "  -4.28%  "
Set-TextValue $ws.Cells.Item(24, 4) "22.013.05"
$ws.Cells.Item(24, 5).Value = "  -1.57%  "
Set-TextValue $ws.Cells.Item(25, 4) "2.357"
$ws.Cells.Item(25, 5).Value = "  +2.51%  "
Set-TextValue $ws.Cells.Item(26, 4) "2.408"
$ws.Cells.Item(26, 5).Value = "  -4.56%  "
Set-TextValue $ws.Cells.Item(27, 4) "147.54"
$ws.Cells.Item(27, 5).Value = "  -2.13%  "
Set-TextValue $ws.Cells.Item(28, 4) "18.50"
$ws.Cells.Item(28, 5).Value = "  -4.76%  "
Set-TextValue $ws.Cells.Item(29, 4) "4.834"
$ws.Cells.Item(29, 5).Value = "  -1.38%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.722.58"
$ws.Cells.Item(30, 5).Value = "  -0.60%  "
Set-TextValue $ws.Cells.Item(31, 4) "118.51"
$ws.Cells.Item(31, 5).Value = "  -3.87%  "
Set-TextValue $ws.Cells.Item(32, 4) "1.045"
$ws.Cells.Item(32, 5).Value = "  +0.03%  "
Set-TextValue $ws.Cells.Item(33, 4) "5.658"
$ws.Cells.Item(33, 5).Value = "  -4.53%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.08343"
Set-TextValue $ws.Cells.Item(35, 4) "9.143"
$ws.Cells.Item(35, 5).Value = "  -3.65%  "
Set-TextValue $ws.Cells.Item(36, 4) "1.596"
$ws.Cells.Item(36, 5).Value = "  -16.55%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.06142"
$ws.Cells.Item(37, 5).Value = "  -2.54%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.02258"
$ws.Cells.Item(38, 5).Value = "  -5.07%  "
Set-TextValue $ws.Cells.Item(39, 4) "5.085"
$ws.Cells.Item(39, 5).Value = "  -3.47%  "
Set-TextValue $ws.Cells.Item(40, 4) "1.212"
$ws.Cells.Item(40, 5).Value = "  -4.67%  "
Set-TextValue $ws.Cells.Item(41, 4) "0.2061"
$ws.Cells.Item(41, 5).Value = "  -5.57%  "
$ws.Cells.Item(42, 5).Value = "  -0.02%  "
$ws.Cells.Item(43, 5).Value = "  -4.31%  "
Set-TextValue $ws.Cells.Item(44, 4) "0.5788"
$ws.Cells.Item(44, 5).Value = "  -4.82%  "
Set-TextValue $ws.Cells.Item(45, 4) "13.11"
$ws.Cells.Item(45, 5).Value = "  -4.24%  "
Set-TextValue $ws.Cells.Item(46, 4) "3.710"
$ws.Cells.Item(46, 5).Value = "  -1.03%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.5545"
$ws.Cells.Item(47, 5).Value = "  -5.90%  "
Set-TextValue $ws.Cells.Item(48, 4) "117.78"
$ws.Cells.Item(48, 5).Value = "  -5.04%  "
Set-TextValue $ws.Cells.Item(49, 4) "1.887"
$ws.Cells.Item(49, 5).Value = "  -5.69%  "
Set-TextValue $ws.Cells.Item(50, 4) "1.134"
$ws.Cells.Item(50, 5).Value = "  -3.85%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.06815"
$ws.Cells.Item(51, 5).Value = "  -4.44%  "
